# Auto-generated edit script: updates currentAveragePrice-derived
# columns (H-N) across multiple job sheets per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 457.06668
$ws.Range("I9").Value = 307.25
$ws.Range("J9").Value = 1056.3334
$ws.Range("K9").Value = 307.25
$ws.Range("L9").Value = 1056.3334
$ws.Range("M9").Value = -138.25
$ws.Range("N9").Value = -1394.3334
$ws.Range("H19").Value = 1483.1875
$ws.Range("I19").Value = 1401.75
$ws.Range("J19").Value = 1510.3334
$ws.Range("K19").Value = 1401.75
$ws.Range("L19").Value = 1510.3334
$ws.Range("M19").Value = -1226.75
$ws.Range("N19").Value = -1860.3334
$ws.Range("H41").Value = 601
$ws.Range("I41").Value = 471
$ws.Range("K41").Value = 471
$ws.Range("M41").Value = -31
$ws.Range("H121").Value = 2480.889
$ws.Range("J121").Value = 2697.875
$ws.Range("L121").Value = 8093.625
$ws.Range("N121").Value = -11587.625
$ws.Range("H127").Value = 1007.0833
$ws.Range("I127").Value = 644.0909
$ws.Range("K127").Value = 1932.2727
$ws.Range("M127").Value = 3027.7273
$ws.Range("H132").Value = 2744.5134
$ws.Range("I132").Value = 2404.111
$ws.Range("K132").Value = 7212.333
$ws.Range("M132").Value = -4682.333
$ws.Range("H137").Value = 1745.5161
$ws.Range("I137").Value = 1246.7646
$ws.Range("K137").Value = 3740.2938
$ws.Range("M137").Value = -1190.2938
$ws.Range("H141").Value = 3766.0222
$ws.Range("I141").Value = 2319.5557
$ws.Range("K141").Value = 6958.6671
$ws.Range("M141").Value = -1778.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23187.424
$ws.Range("I32").Value = 11730.737
$ws.Range("K32").Value = 11730.737
$ws.Range("M32").Value = -11443.737
$ws.Range("H74").Value = 12151.125
$ws.Range("I74").Value = 4004
$ws.Range("J74").Value = 17039.4
$ws.Range("K74").Value = 4004
$ws.Range("L74").Value = 17039.4
$ws.Range("M74").Value = -3130
$ws.Range("N74").Value = -18787.4
$ws.Range("H77").Value = 12151.125
$ws.Range("I77").Value = 4004
$ws.Range("J77").Value = 17039.4
$ws.Range("K77").Value = 20020
$ws.Range("L77").Value = 85197
$ws.Range("M77").Value = -15652
$ws.Range("N77").Value = -93933
$ws.Range("H102").Value = 11749.75
$ws.Range("I102").Value = 11749.75
$ws.Range("K102").Value = 11749.75
$ws.Range("M102").Value = -10127.75
$ws.Range("H105").Value = 25000
$ws.Range("J105").Value = 25000
$ws.Range("L105").Value = 25000
$ws.Range("N105").Value = -31988
$ws.Range("H110").Value = 1334.1538
$ws.Range("I110").Value = 1362
$ws.Range("K110").Value = 1362
$ws.Range("M110").Value = 683
$ws.Range("H132").Value = 6316.8
$ws.Range("I132").Value = 4952.5713
$ws.Range("K132").Value = 14857.7139
$ws.Range("M132").Value = -12327.7139
$ws.Range("H133").Value = 80407
$ws.Range("J133").Value = 80407
$ws.Range("L133").Value = 80407
$ws.Range("N133").Value = -85467

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1643
$ws.Range("I20").Value = 1592
$ws.Range("K20").Value = 1592
$ws.Range("M20").Value = -1345
$ws.Range("H59").Value = 92554.664
$ws.Range("J59").Value = 94124
$ws.Range("L59").Value = 94124
$ws.Range("N59").Value = -95818
$ws.Range("H132").Value = 89898
$ws.Range("J132").Value = 89898
$ws.Range("L132").Value = 89898
$ws.Range("N132").Value = -100018
$ws.Range("H134").Value = 8280.883
$ws.Range("I134").Value = 5349.75
$ws.Range("J134").Value = 10886.333
$ws.Range("K134").Value = 16049.25
$ws.Range("L134").Value = 32658.999
$ws.Range("M134").Value = -13514.25
$ws.Range("N134").Value = -37728.999
$ws.Range("H140").Value = 113999.5
$ws.Range("J140").Value = 113999.5
$ws.Range("L140").Value = 113999.5
$ws.Range("N140").Value = -124359.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 77974.75
$ws.Range("J52").Value = 83999.5
$ws.Range("L52").Value = 83999.5
$ws.Range("N52").Value = -84587.5
$ws.Range("H105").Value = 1291.875
$ws.Range("I105").Value = 1376.7142
$ws.Range("J105").Value = 698
$ws.Range("K105").Value = 1376.7142
$ws.Range("L105").Value = 698
$ws.Range("M105").Value = 370.2858000000001
$ws.Range("N105").Value = -4192
$ws.Range("H135").Value = 82153
$ws.Range("J135").Value = 82153
$ws.Range("L135").Value = 82153
$ws.Range("N135").Value = -92293

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 376.2
$ws.Range("I2").Value = 703.1667
$ws.Range("J2").Value = 158.22223
$ws.Range("K2").Value = 4219.0002
$ws.Range("L2").Value = 949.33338
$ws.Range("M2").Value = -4106.0002
$ws.Range("N2").Value = -1175.33338
$ws.Range("H5").Value = 1345.6471
$ws.Range("I5").Value = 1468.1
$ws.Range("K5").Value = 4404.299999999999
$ws.Range("M5").Value = -4292.299999999999
$ws.Range("H20").Value = 2996
$ws.Range("J20").Value = 2996
$ws.Range("L20").Value = 8988
$ws.Range("N20").Value = -9442
$ws.Range("H22").Value = 1279.8
$ws.Range("J22").Value = 833
$ws.Range("L22").Value = 2499
$ws.Range("N22").Value = -2837
$ws.Range("H27").Value = 1279.8
$ws.Range("J27").Value = 833
$ws.Range("L27").Value = 2499
$ws.Range("N27").Value = -2703
$ws.Range("H33").Value = 21.75
$ws.Range("I33").Value = 7.3333335
$ws.Range("K33").Value = 44.000001
$ws.Range("M33").Value = 238.999999
$ws.Range("H34").Value = 8763.846
$ws.Range("J34").Value = 10075.546
$ws.Range("L34").Value = 30226.638
$ws.Range("N34").Value = -30394.638
$ws.Range("H46").Value = 2021
$ws.Range("I46").Value = 1016.5
$ws.Range("J46").Value = 2882
$ws.Range("K46").Value = 3049.5
$ws.Range("L46").Value = 8646
$ws.Range("M46").Value = -2958.5
$ws.Range("N46").Value = -8828
$ws.Range("H51").Value = 2250
$ws.Range("I51").Value = 1500
$ws.Range("K51").Value = 4500
$ws.Range("M51").Value = -4040
$ws.Range("H104").Value = 789
$ws.Range("J104").Value = 789
$ws.Range("L104").Value = 2367
$ws.Range("N104").Value = -7609
$ws.Range("H131").Value = 4081.6924
$ws.Range("I131").Value = 4510
$ws.Range("K131").Value = 13530
$ws.Range("M131").Value = -8490
$ws.Range("H135").Value = 1345.6471
$ws.Range("I135").Value = 1468.1
$ws.Range("K135").Value = 13212.9
$ws.Range("M135").Value = -10677.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9123.714
$ws.Range("I70").Value = 8469.75
$ws.Range("J70").Value = 9995.666999999999
$ws.Range("K70").Value = 8469.75
$ws.Range("L70").Value = 9995.666999999999
$ws.Range("M70").Value = -8199.75
$ws.Range("N70").Value = -10535.667
$ws.Range("H73").Value = 9123.714
$ws.Range("I73").Value = 8469.75
$ws.Range("J73").Value = 9995.666999999999
$ws.Range("K73").Value = 8469.75
$ws.Range("L73").Value = 9995.666999999999
$ws.Range("M73").Value = -7533.75
$ws.Range("N73").Value = -11867.667
$ws.Range("H102").Value = 3432.4614
$ws.Range("I102").Value = 1764
$ws.Range("K102").Value = 1764
$ws.Range("M102").Value = -142
$ws.Range("H126").Value = 6829.1035
$ws.Range("I126").Value = 6900.6665
$ws.Range("J126").Value = 6752.4287
$ws.Range("K126").Value = 20701.9995
$ws.Range("L126").Value = 20257.2861
$ws.Range("M126").Value = -18231.9995
$ws.Range("N126").Value = -25197.2861
$ws.Range("H132").Value = 7333.778
$ws.Range("I132").Value = 6746.3335
$ws.Range("K132").Value = 20239.0005
$ws.Range("M132").Value = -17709.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5315.448
$ws.Range("I7").Value = 5153.6294
$ws.Range("K7").Value = 5153.6294
$ws.Range("M7").Value = -5041.6294
$ws.Range("H55").Value = 1913.6
$ws.Range("I55").Value = 516.75
$ws.Range("J55").Value = 2844.8333
$ws.Range("K55").Value = 516.75
$ws.Range("L55").Value = 2844.8333
$ws.Range("M55").Value = -343.75
$ws.Range("N55").Value = -3190.8333
$ws.Range("H61").Value = 3084.257
$ws.Range("I61").Value = 2457.2083
$ws.Range("J61").Value = 4452.364
$ws.Range("K61").Value = 2457.2083
$ws.Range("L61").Value = 4452.364
$ws.Range("M61").Value = -2255.2083
$ws.Range("N61").Value = -4856.364
$ws.Range("H113").Value = 3084.257
$ws.Range("I113").Value = 2457.2083
$ws.Range("J113").Value = 4452.364
$ws.Range("K113").Value = 2457.2083
$ws.Range("L113").Value = 4452.364
$ws.Range("M113").Value = -287.2082999999998
$ws.Range("N113").Value = -8792.364
$ws.Range("H126").Value = 5315.448
$ws.Range("I126").Value = 5153.6294
$ws.Range("K126").Value = 15460.8882
$ws.Range("M126").Value = -12990.8882
$ws.Range("H136").Value = 7021.7124
$ws.Range("I136").Value = 2881.5833
$ws.Range("K136").Value = 8644.749899999999
$ws.Range("M136").Value = -6094.749899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 29772.027
$ws.Range("I81").Value = 57033.11
$ws.Range("J81").Value = 3945.7368
$ws.Range("K81").Value = 114066.22
$ws.Range("L81").Value = 7891.4736
$ws.Range("M81").Value = -113005.22
$ws.Range("N81").Value = -10013.4736
$ws.Range("H84").Value = 29772.027
$ws.Range("I84").Value = 57033.11
$ws.Range("J84").Value = 3945.7368
$ws.Range("K84").Value = 570331.1
$ws.Range("L84").Value = 39457.368
$ws.Range("M84").Value = -565027.1
$ws.Range("N84").Value = -50065.368
$ws.Range("H132").Value = 4062.7727
$ws.Range("I132").Value = 2633.2666
$ws.Range("J132").Value = 7126
$ws.Range("K132").Value = 7899.7998
$ws.Range("L132").Value = 21378
$ws.Range("M132").Value = -5369.7998
$ws.Range("N132").Value = -26438

